$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8 (ALC)
$ws.Cells.Item(8, 8).Value = 1388.3636
$ws.Cells.Item(8, 9).Value = 474.22223
$ws.Cells.Item(8, 10).Value = 5502
$ws.Cells.Item(8, 11).Value = 1422.66669
$ws.Cells.Item(8, 12).Value = 16506
$ws.Cells.Item(8, 13).Value = -1283.66669
$ws.Cells.Item(8, 14).Value = -16784

# Row 9 (ALC)
$ws.Cells.Item(9, 8).Value = 98.666664
$ws.Cells.Item(9, 9).Value = 156
$ws.Cells.Item(9, 10).Value = 52.8
$ws.Cells.Item(9, 11).Value = 156
$ws.Cells.Item(9, 12).Value = 52.8
$ws.Cells.Item(9, 13).Value = 13
$ws.Cells.Item(9, 14).Value = -390.8

# Row 31 (ALC)
$ws.Cells.Item(31, 8).Value = 300
$ws.Cells.Item(31, 9).Value = 300
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 900
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).ClearContents()
$ws.Cells.Item(31, 14).ClearContents()

# Row 53 (ALC)
$ws.Cells.Item(53, 8).Value = 178.27272
$ws.Cells.Item(53, 9).Value = 186.83333
$ws.Cells.Item(53, 10).Value = 168
$ws.Cells.Item(53, 11).Value = 186.83333
$ws.Cells.Item(53, 12).Value = 168
$ws.Cells.Item(53, 13).Value = 450.16667
$ws.Cells.Item(53, 14).Value = -1442

# Row 112 (ALC)
$ws.Cells.Item(112, 8).Value = 5587.0586
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 5587.0586
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 16761.1758
$ws.Cells.Item(112, 14).Value = -18977.1758

# Row 129 (ALC)
$ws.Cells.Item(129, 8).Value = 1298.65
$ws.Cells.Item(129, 9).Value = 575.7692
$ws.Cells.Item(129, 10).Value = 1646.7037
$ws.Cells.Item(129, 11).Value = 1727.3076
$ws.Cells.Item(129, 12).Value = 4940.1111
$ws.Cells.Item(129, 13).Value = 3272.6924
$ws.Cells.Item(129, 14).Value = -14940.1111

# Row 135 (ALC)
$ws.Cells.Item(135, 8).Value = 1139.3572
$ws.Cells.Item(135, 9).Value = 1124.9
$ws.Cells.Item(135, 10).Value = 1175.5
$ws.Cells.Item(135, 11).Value = 10124.1
$ws.Cells.Item(135, 12).Value = 10579.5
$ws.Cells.Item(135, 13).Value = -7589.1
$ws.Cells.Item(135, 14).Value = -15649.5

# Row 137 (ALC)
$ws.Cells.Item(137, 8).Value = 1595.5217
$ws.Cells.Item(137, 9).Value = 1680.05
$ws.Cells.Item(137, 10).Value = 1530.5
$ws.Cells.Item(137, 11).Value = 5040.15
$ws.Cells.Item(137, 12).Value = 4591.5
$ws.Cells.Item(137, 13).Value = -2490.15
$ws.Cells.Item(137, 14).Value = -9691.5

$ws = $wb.Worksheets.Item("ARM")
# Row 110 (ARM)
$ws.Cells.Item(110, 8).Value = 151049.67
$ws.Cells.Item(110, 9).Value = 151049.67
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 151049.67
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = -149004.67

# Row 122 (ARM)
$ws.Cells.Item(122, 8).Value = 1793.9166
$ws.Cells.Item(122, 9).Value = 1467
$ws.Cells.Item(122, 10).Value = 3428.5
$ws.Cells.Item(122, 11).Value = 4401
$ws.Cells.Item(122, 12).Value = 10285.5
$ws.Cells.Item(122, 13).Value = -1951
$ws.Cells.Item(122, 14).Value = -15185.5

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (BSM)
$ws.Cells.Item(107, 8).Value = 1232.3334
$ws.Cells.Item(107, 9).Value = 1023.5
$ws.Cells.Item(107, 10).Value = 1650
$ws.Cells.Item(107, 11).Value = 1023.5
$ws.Cells.Item(107, 12).Value = 1650
$ws.Cells.Item(107, 13).Value = 896.5
$ws.Cells.Item(107, 14).Value = -5490

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Cells.Item(16, 8).Value = 1570.25
$ws.Cells.Item(16, 9).Value = 1530
$ws.Cells.Item(16, 10).Value = 1599
$ws.Cells.Item(16, 11).Value = 1530
$ws.Cells.Item(16, 12).Value = 1599
$ws.Cells.Item(16, 13).Value = -1243
$ws.Cells.Item(16, 14).Value = -2173

# Row 31 (CRP)
$ws.Cells.Item(31, 8).Value = 4351.4854
$ws.Cells.Item(31, 9).Value = 1373.7576
$ws.Cells.Item(31, 10).Value = 7159.057
$ws.Cells.Item(31, 11).Value = 1373.7576
$ws.Cells.Item(31, 12).Value = 7159.057
$ws.Cells.Item(31, 13).Value = -1078.7576
$ws.Cells.Item(31, 14).Value = -7749.057

# Row 34 (CRP)
$ws.Cells.Item(34, 8).Value = 4351.4854
$ws.Cells.Item(34, 9).Value = 1373.7576
$ws.Cells.Item(34, 10).Value = 7159.057
$ws.Cells.Item(34, 11).Value = 1373.7576
$ws.Cells.Item(34, 12).Value = 7159.057
$ws.Cells.Item(34, 13).Value = -1171.7576
$ws.Cells.Item(34, 14).Value = -7563.057

# Row 58 (CRP)
$ws.Cells.Item(58, 8).Value = 1123.125
$ws.Cells.Item(58, 9).Value = 829.875
$ws.Cells.Item(58, 10).Value = 1856.25
$ws.Cells.Item(58, 11).Value = 829.875
$ws.Cells.Item(58, 12).Value = 1856.25
$ws.Cells.Item(58, 13).Value = -626.875
$ws.Cells.Item(58, 14).Value = -2262.25

# Row 107 (CRP)
$ws.Cells.Item(107, 8).Value = 2500778.5
$ws.Cells.Item(107, 9).Value = 5682297.5
$ws.Cells.Item(107, 10).Value = 1013.7143
$ws.Cells.Item(107, 11).Value = 5682297.5
$ws.Cells.Item(107, 12).Value = 1013.7143
$ws.Cells.Item(107, 13).Value = -5680377.5
$ws.Cells.Item(107, 14).Value = -4853.7143

# Row 113 (CRP)
$ws.Cells.Item(113, 8).Value = 1570.25
$ws.Cells.Item(113, 9).Value = 1530
$ws.Cells.Item(113, 10).Value = 1599
$ws.Cells.Item(113, 11).Value = 1530
$ws.Cells.Item(113, 12).Value = 1599
$ws.Cells.Item(113, 13).Value = 640
$ws.Cells.Item(113, 14).Value = -5939

# Row 132 (CRP)
$ws.Cells.Item(132, 8).Value = 6668798.5
$ws.Cells.Item(132, 9).Value = 1952.3846
$ws.Cells.Item(132, 10).Value = 13891215
$ws.Cells.Item(132, 11).Value = 5857.1538
$ws.Cells.Item(132, 12).Value = 41673645
$ws.Cells.Item(132, 13).Value = -3327.1538
$ws.Cells.Item(132, 14).Value = -41678705

# Row 136 (CRP)
$ws.Cells.Item(136, 8).Value = 1123.125
$ws.Cells.Item(136, 9).Value = 829.875
$ws.Cells.Item(136, 10).Value = 1856.25
$ws.Cells.Item(136, 11).Value = 2489.625
$ws.Cells.Item(136, 12).Value = 5568.75
$ws.Cells.Item(136, 13).Value = 60.375
$ws.Cells.Item(136, 14).Value = -10668.75

$ws = $wb.Worksheets.Item("GSM")
# Row 38 (GSM)
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 12).ClearContents()
$ws.Cells.Item(38, 14).ClearContents()

# Row 40 (GSM)
$ws.Cells.Item(40, 8).Value = 12500
$ws.Cells.Item(40, 9).Value = 20000
$ws.Cells.Item(40, 10).Value = 5000
$ws.Cells.Item(40, 11).Value = 20000
$ws.Cells.Item(40, 12).Value = 5000
$ws.Cells.Item(40, 13).Value = -19849
$ws.Cells.Item(40, 14).Value = -5302

# Row 44 (GSM)
$ws.Cells.Item(44, 8).Value = 6000
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 6000
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).ClearContents()
$ws.Cells.Item(44, 13).ClearContents()
$ws.Cells.Item(44, 14).Value = -7192

# Row 57 (GSM)
$ws.Cells.Item(57, 8).Value = 18999.834
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 18999.834
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 18999.834
$ws.Cells.Item(57, 14).Value = -20639.834

# Row 58 (GSM)
$ws.Cells.Item(58, 8).Value = 16000
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 16000
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 12).Value = 16000
$ws.Cells.Item(58, 14).Value = -16554

# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 2909.1052
$ws.Cells.Item(132, 9).Value = 2326.889
$ws.Cells.Item(132, 10).Value = 4338.1816
$ws.Cells.Item(132, 11).Value = 6980.667
$ws.Cells.Item(132, 12).Value = 13014.5448
$ws.Cells.Item(132, 13).Value = -4450.667
$ws.Cells.Item(132, 14).Value = -18074.5448

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Cells.Item(61, 8).Value = 3024.625
$ws.Cells.Item(61, 9).Value = 2801.8333
$ws.Cells.Item(61, 10).Value = 3693
$ws.Cells.Item(61, 11).Value = 2801.8333
$ws.Cells.Item(61, 12).Value = 3693
$ws.Cells.Item(61, 13).Value = -2599.8333
$ws.Cells.Item(61, 14).Value = -4097

# Row 113 (LTW)
$ws.Cells.Item(113, 8).Value = 3024.625
$ws.Cells.Item(113, 9).Value = 2801.8333
$ws.Cells.Item(113, 10).Value = 3693
$ws.Cells.Item(113, 11).Value = 2801.8333
$ws.Cells.Item(113, 12).Value = 3693
$ws.Cells.Item(113, 13).Value = -631.8332999999998
$ws.Cells.Item(113, 14).Value = -8033

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (WVR)
$ws.Cells.Item(2, 8).Value = 167501170
$ws.Cells.Item(2, 9).Value = 5000000
$ws.Cells.Item(2, 10).Value = 200001390
$ws.Cells.Item(2, 11).Value = 5000000
$ws.Cells.Item(2, 12).Value = 200001390
$ws.Cells.Item(2, 13).Value = -4999888
$ws.Cells.Item(2, 14).Value = -200001614

# Row 9 (WVR)
$ws.Cells.Item(9, 8).Value = 56064.2
$ws.Cells.Item(9, 9).Value = 300
$ws.Cells.Item(9, 10).Value = 70005.25
$ws.Cells.Item(9, 11).Value = 300
$ws.Cells.Item(9, 12).Value = 70005.25
$ws.Cells.Item(9, 13).Value = -160
$ws.Cells.Item(9, 14).Value = -70285.25

# Row 20 (WVR)
$ws.Cells.Item(20, 8).Value = 44847
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 44847
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 44847
$ws.Cells.Item(20, 14).Value = -45327

# Row 26 (WVR)
$ws.Cells.Item(26, 8).Value = 70014
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 70014
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 70014
$ws.Cells.Item(26, 14).Value = -70600

# Row 34 (WVR)
$ws.Cells.Item(34, 8).Value = 15000
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 15000
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).ClearContents()
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(34, 14).Value = -15406

# Row 39 (WVR)
$ws.Cells.Item(39, 8).Value = 56715.668
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 56715.668
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 56715.668
$ws.Cells.Item(39, 14).Value = -57541.668

# Row 40 (WVR)
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).ClearContents()
$ws.Cells.Item(40, 14).ClearContents()

# Row 51 (WVR)
$ws.Cells.Item(51, 8).Value = 13842.857
$ws.Cells.Item(51, 9).Value = 7900
$ws.Cells.Item(51, 10).Value = 14833.333
$ws.Cells.Item(51, 11).Value = 7900
$ws.Cells.Item(51, 12).Value = 14833.333
$ws.Cells.Item(51, 13).Value = -7390
$ws.Cells.Item(51, 14).Value = -15853.333

# Row 132 (WVR)
$ws.Cells.Item(132, 8).Value = 5720863.5
$ws.Cells.Item(132, 9).Value = 2169.4333
$ws.Cells.Item(132, 10).Value = 13890427
$ws.Cells.Item(132, 11).Value = 6508.2999
$ws.Cells.Item(132, 12).Value = 41671281
$ws.Cells.Item(132, 13).Value = -3978.2999
$ws.Cells.Item(132, 14).Value = -41676341

# Row 136 (WVR)
$ws.Cells.Item(136, 8).Value = 15154622
$ws.Cells.Item(136, 9).Value = 10872812
$ws.Cells.Item(136, 10).Value = 25002788
$ws.Cells.Item(136, 11).Value = 32618436
$ws.Cells.Item(136, 12).Value = 75008364
$ws.Cells.Item(136, 13).Value = -32615886
$ws.Cells.Item(136, 14).Value = -75013464
